# Insert a new header row above the existing data (row 1), shifting
# all existing rows down by one, then populate the new header cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "revenue"
$ws.Range("C1").Value = "budget"
$ws.Range("D1").Value = "roi"

$ws.Range("I14").Select()
